# Renumber vertices P136..P260 -> P000..P124 (zero-padded to 3 digits)
# throughout the body paragraph that describes the perimeter, while
# leaving the document title ("Sede Urbana P136 até P260") untouched.

$d = $word.ActiveDocument

# The title lives in paragraph 1 (Heading1 style) and must stay as-is.
# All of the vertex descriptions live in paragraph 2 - scope every
# Find/Replace to that paragraph's Range so the heading is never touched.

for ($n = 136; $n -le 260; $n++) {
    $newNum  = $n - 136
    $oldText = "P" + $n
    $newText = "P" + $newNum.ToString().PadLeft(3, "0")

    $rng = $d.Paragraphs.Item(2).Range
    $rng.Find.ClearFormatting()
    [void]$rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
}

Write-Host "Renumbering complete"
